# Update Name of Algo
# Apply value corrections to the KNN imputation result data on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 12.802
$ws.Range("D9").Value = -7.946000000000001
$ws.Range("D18").Value = -8.144
$ws.Range("D20").Value = -8.044
